# The workbook gained two new price-record rows (Murcott / 18kg tray,
# "Región de O'Higgins") that were inserted right before the existing
# row 353 ("Fecha" 2021-08-06 / Clementina). Everything that used to be
# row 353 onward (through the old row 441) shifts down by two rows, and
# the sheet's used range grows from A1:T441 to A1:T443.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the existing row 353, pushing old rows
# 353..441 down to 355..443.
$ws.Rows("353:354").Insert()

# --- New row 353 -----------------------------------------------------
$ws.Range("A353").Value = 7
$ws.Range("B353").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C353").Value = "Ñuble"
$ws.Range("D353").Value = 45204
$ws.Range("E353").Value = 16
$ws.Range("F353").Value = "Fruta"
$ws.Range("G353").Value = 100102
$ws.Range("H353").Value = "Cítricos"
$ws.Range("I353").Value = 100102004
$ws.Range("J353").Value = "Mandarina"
$ws.Range("K353").Value = "Murcott"
$ws.Range("L353").Value = "Primera"
$ws.Range("M353").Value = 120
$ws.Range("N353").Value = 9000
$ws.Range("O353").Value = 9000
$ws.Range("P353").Value = 9000
$ws.Range("Q353").Value = "`$/bandeja 18 kilos"
$ws.Range("R353").Value = "Región de O'Higgins"
$ws.Range("S353").Value = 500
$ws.Range("T353").Value = 18

# --- New row 354 -----------------------------------------------------
$ws.Range("A354").Value = 7
$ws.Range("B354").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C354").Value = "Ñuble"
$ws.Range("D354").Value = 45204
$ws.Range("E354").Value = 16
$ws.Range("F354").Value = "Fruta"
$ws.Range("G354").Value = 100102
$ws.Range("H354").Value = "Cítricos"
$ws.Range("I354").Value = 100102004
$ws.Range("J354").Value = "Mandarina"
$ws.Range("K354").Value = "Murcott"
$ws.Range("L354").Value = "Segunda"
$ws.Range("M354").Value = 120
$ws.Range("N354").Value = 7000
$ws.Range("O354").Value = 7000
$ws.Range("P354").Value = 7000
$ws.Range("Q354").Value = "`$/bandeja 18 kilos"
$ws.Range("R354").Value = "Región de O'Higgins"
$ws.Range("S354").Value = 389
$ws.Range("T354").Value = 18

# Keep the date columns formatted like the rest of the "Fecha" column.
$ws.Range("D353:D354").NumberFormat = $ws.Range("D352").NumberFormat
